# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '55.857.18'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -3.17%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.925.29'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -3.55%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '504.29'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.73%  '
$ws.Range('E6').Value = '  -4.91%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.423'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -4.34%  '
$ws.Range('E9').Value = '  -4.94%  '
$ws.Range('E10').Value = '  -5.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.352'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.51%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '3.426.79'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.56%  '
$ws.Range('E13').Value = '  -3.99%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '25.91'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.08%  '
$ws.Range('E15').Value = '  -3.33%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '55.733.25'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.50%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.98'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -4.83%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.922.84'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -3.74%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.72'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.14%  '
$ws.Range('E20').Value = '  -3.93%  '
$ws.Range('E21').Value = '  -5.48%  '
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.488'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.47%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '62.96'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.87%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.045.71'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.51%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').Value = '  -4.79%  '
$ws.Range('E28').Value = '  -9.74%  '
$ws.Range('E29').Value = '  -6.98%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.85'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -8.93%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.77'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.69%  '
$ws.Range('E32').Value = '  -5.98%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.89'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.45%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '150.84'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.19%  '
$ws.Range('E35').Value = '  -7.21%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.63'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.50%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '24.13'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.82%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.20'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -7.03%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0646'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -5.73%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '36.41'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.95%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.71'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -4.66%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.640'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.52%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.124.99'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -8.33%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '6.00'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.34'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -6.38%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.925'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -7.17%  '
$ws.Range('E48').Value = '  -2.32%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '18.72'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.69%  '
$ws.Range('E50').Value = '  -6.42%  '
$ws.Range('E51').Value = '  -9.71%  '
